# Auto-generated script applying the cryptos.xlsx diff via Excel COM interop
# (Sun Feb  4 10:24:45 UTC 2024 GitHub Actions crypto-price refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.913.13"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.301.74"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.511"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.118"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").Value = "2.660.75"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "2.303.19"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "42.852.18"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.58%  "
$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  -4.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.41%  "
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.109"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "2.011.82"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0280"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.84%  "
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.526.37"
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.79%  "
